$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New UF / Localidade rows (rows 7-12), entering shared strings in the
# order that matches the original authoring session so the sharedStrings
# table is rebuilt with the same index order.
$ws.Range("A7").Value = "SP"
$ws.Range("B7").Value = "São Paulo"

$ws.Range("A8").Value = "RJ"
$ws.Range("B8").Value = "Rio de Janeiro"

$ws.Range("A9").Value = "RS"
$ws.Range("B9").Value = "Porto Alegre"

$ws.Range("A10").Value = "MS"
$ws.Range("B10").Value = "Campo Grande"

$ws.Range("B11").Value = "Natal"
$ws.Range("B12").Value = "Parnamirim"
$ws.Range("A11").Value = "RN"
$ws.Range("A12").Value = "RN"

# Selection cursor ends on B18 in the saved file.
$ws.Range("B18").Select()
